$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.012.29"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "2.626.06"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'111.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").Value = "'322.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("D7").Value = "'0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.87%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "

$ws.Range("D10").Value = "'39.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "

$ws.Range("D11").Value = "'19.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.83%  "

$ws.Range("D12").Value = "'0.0809"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'7.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "

$ws.Range("D15").Value = "3.040.47"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "2.634.60"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "'0.856"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").Value = "48.990.17"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("E20").Value = "  -3.55%  "

$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "0.0₃0942"
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("D23").Value = "'268.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.38%  "

$ws.Range("D24").Value = "'68.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.95%  "

$ws.Range("D25").Value = "'2.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.66%  "

$ws.Range("D26").Value = "'26.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.44%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "'10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").Value = "'35.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "

$ws.Range("E31").Value = "  -4.96%  "

$ws.Range("D32").Value = "'49.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("D33").Value = "'5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").Value = "'18.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.34%  "

$ws.Range("D37").Value = "'4.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.55%  "

$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").Value = "'3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").Value = "'127.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "

$ws.Range("E41").Value = "  -2.09%  "

$ws.Range("D42").Value = "'22.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.51%  "

$ws.Range("E43").Value = "  -4.27%  "

$ws.Range("D44").Value = "'0.0316"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "2.058.92"
$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("E46").Value = "  +6.51%  "

$ws.Range("D47").Value = "'3.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.65%  "

$ws.Range("E48").Value = "  -4.37%  "

$ws.Range("D49").Value = "'8.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("E51").Value = "  +1.30%  "
